$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.998.47'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.616.00'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +1.43%  '
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '525.20'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.43%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '154.96'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.50%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.589'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +1.24%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '6.65'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.26%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.106'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.20%  '
$ws.Range('E11').Value = '  -0.06%  '
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.072.30'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.01%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '61.029.08'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '21.67'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.06%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000142'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +1.21%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.618.56'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '4.77'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -0.15%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '356.05'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.68%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.59'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.80%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.22'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.68%  '
$ws.Range('E22').Value = '  +0.19%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '61.14'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.80%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.428'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.63%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.167'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.03%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.728.32'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.71%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0851'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.41'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -0.13%  '
$ws.Range('E30').Value = '  -0.07%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.27'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +9.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '19.49'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.45%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.61'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.00%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '150.29'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -1.74%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.16'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.54%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.20'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.75%  '
$ws.Range('B37').Value = 'SuiNetwork'
$ws.Range('C37').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.917'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +7.65%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.906'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +7.15%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.50'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.14%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.79'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.67%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '291.54'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -1.81%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.102'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +1.97%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.630'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.21%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0560'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.19%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.997'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.17%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '5.01'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.76%  '
$ws.Range('B47').Value = 'EnergySwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '19.62'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.79%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0238'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.99%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '10.35'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.48%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '19.21'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.39%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.973.42'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.01%  '
